$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3751.7273
$ws.Range("J17").Value = 3751.7273
$ws.Range("L17").Value = 11255.1819
$ws.Range("N17").Value = -11591.1819
$ws.Range("H21").Value = 17
$ws.Range("I21").Value = 17
$ws.Range("K21").Value = 17
$ws.Range("M21").Value = 451
$ws.Range("H23").Value = 17
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = 17
$ws.Range("M23").Value = 217
$ws.Range("H40").Value = 3656.8333
$ws.Range("I40").Value = 3899
$ws.Range("J40").Value = 3172.5
$ws.Range("K40").Value = 3899
$ws.Range("L40").Value = 3172.5
$ws.Range("M40").Value = -3724
$ws.Range("N40").Value = -3522.5
$ws.Range("H42").Value = 2887.6365
$ws.Range("I42").Value = 1220.3334
$ws.Range("J42").Value = 4888.4
$ws.Range("K42").Value = 3661.0002
$ws.Range("L42").Value = 14665.2
$ws.Range("M42").Value = -3431.0002
$ws.Range("N42").Value = -15125.2
$ws.Range("H86").Value = 2464.2307
$ws.Range("I86").Value = 1599.75
$ws.Range("K86").Value = 1599.75
$ws.Range("M86").Value = -476.75
$ws.Range("H89").Value = 2464.2307
$ws.Range("I89").Value = 1599.75
$ws.Range("K89").Value = 7998.75
$ws.Range("M89").Value = -2382.75
$ws.Range("H123").Value = 90000
$ws.Range("J123").Value = 90000
$ws.Range("L123").Value = 90000
$ws.Range("N123").Value = -99800
$ws.Range("H137").Value = 1400.8889
$ws.Range("I137").Value = 1364.84
$ws.Range("K137").Value = 4094.52
$ws.Range("M137").Value = -1544.52
$ws.Range("H141").Value = 2476.6
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 343
$ws.Range("I4").Value = 283.625
$ws.Range("J4").Value = 461.75
$ws.Range("K4").Value = 283.625
$ws.Range("L4").Value = 461.75
$ws.Range("M4").Value = -167.625
$ws.Range("N4").Value = -693.75
$ws.Range("H61").Value = 4060.9534
$ws.Range("I61").Value = 2892.25
$ws.Range("K61").Value = 2892.25
$ws.Range("M61").Value = -2680.25
$ws.Range("H74").Value = 2220.25
$ws.Range("I74").Value = 2167.9614
$ws.Range("K74").Value = 2167.9614
$ws.Range("M74").Value = -1293.9614
$ws.Range("H77").Value = 2220.25
$ws.Range("I77").Value = 2167.9614
$ws.Range("K77").Value = 10839.807
$ws.Range("M77").Value = -6471.807000000001
$ws.Range("H132").Value = 5249.25
$ws.Range("I132").Value = 4332.3335
$ws.Range("K132").Value = 12997.0005
$ws.Range("M132").Value = -10467.0005
$ws.Range("H136").Value = 4060.9534
$ws.Range("I136").Value = 2892.25
$ws.Range("K136").Value = 8676.75
$ws.Range("M136").Value = -6126.75

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3757
$ws.Range("I20").Value = 3660
$ws.Range("K20").Value = 3660
$ws.Range("M20").Value = -3413
$ws.Range("H60").Value = 29780
$ws.Range("J60").Value = 29780
$ws.Range("L60").Value = 29780
$ws.Range("N60").Value = -30978

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6104.826
$ws.Range("I31").Value = 9355.5
$ws.Range("J31").Value = 5201.8613
$ws.Range("K31").Value = 9355.5
$ws.Range("L31").Value = 5201.8613
$ws.Range("M31").Value = -9060.5
$ws.Range("N31").Value = -5791.8613
$ws.Range("H34").Value = 6104.826
$ws.Range("I34").Value = 9355.5
$ws.Range("J34").Value = 5201.8613
$ws.Range("K34").Value = 9355.5
$ws.Range("L34").Value = 5201.8613
$ws.Range("M34").Value = -9153.5
$ws.Range("N34").Value = -5605.8613
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H88").Value = 31982.428
$ws.Range("I88").Value = 14499
$ws.Range("J88").Value = 34896.332
$ws.Range("K88").Value = 14499
$ws.Range("L88").Value = 34896.332
$ws.Range("M88").Value = -14093
$ws.Range("N88").Value = -35708.332
$ws.Range("H91").Value = 31982.428
$ws.Range("I91").Value = 14499
$ws.Range("J91").Value = 34896.332
$ws.Range("K91").Value = 14499
$ws.Range("L91").Value = 34896.332
$ws.Range("M91").Value = -13095
$ws.Range("N91").Value = -37704.332
$ws.Range("H94").Value = 1828.4286
$ws.Range("I94").Value = 1402.3334
$ws.Range("J94").Value = 2148
$ws.Range("K94").Value = 1402.3334
$ws.Range("L94").Value = 2148
$ws.Range("M94").Value = -951.3334
$ws.Range("N94").Value = -3050
$ws.Range("H134").Value = 812.3570999999999
$ws.Range("I134").Value = 734.7826
$ws.Range("J134").Value = 1169.2
$ws.Range("K134").Value = 2204.3478
$ws.Range("L134").Value = 3507.6
$ws.Range("M134").Value = 330.6522
$ws.Range("N134").Value = -8577.6
$ws.Range("H141").Value = 372666.66
$ws.Range("J141").Value = 372666.66
$ws.Range("L141").Value = 372666.66
$ws.Range("N141").Value = -383026.66

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 4607.6665
$ws.Range("J93").Value = 6499.5
$ws.Range("L93").Value = 19498.5
$ws.Range("N93").Value = -23242.5
$ws.Range("H104").Value = 333
$ws.Range("I104").Value = 333
$ws.Range("K104").Value = 999
$ws.Range("M104").Value = 1622
$ws.Range("H107").Value = 375.45456
$ws.Range("I107").Value = 245
$ws.Range("J107").Value = 424.375
$ws.Range("K107").Value = 735
$ws.Range("L107").Value = 1273.125
$ws.Range("M107").Value = 1185
$ws.Range("N107").Value = -5113.125
$ws.Range("H109").Value = 4006.3333
$ws.Range("I109").Value = 3579.5715
$ws.Range("K109").Value = 10738.7145
$ws.Range("M109").Value = -9698.7145
$ws.Range("H118").Value = 6335.364
$ws.Range("I118").Value = 2939.6
$ws.Range("J118").Value = 9165.166999999999
$ws.Range("K118").Value = 8818.799999999999
$ws.Range("L118").Value = 27495.501
$ws.Range("M118").Value = -7575.799999999999
$ws.Range("N118").Value = -29981.501
$ws.Range("H134").Value = 9370.352999999999
$ws.Range("I134").Value = 7966.6
$ws.Range("J134").Value = 19898.5
$ws.Range("K134").Value = 23899.8
$ws.Range("L134").Value = 59695.5
$ws.Range("M134").Value = -18829.8
$ws.Range("N134").Value = -69835.5
$ws.Range("H137").Value = 5458.7144
$ws.Range("I137").Value = 5091.4
$ws.Range("J137").Value = 5662.778
$ws.Range("K137").Value = 15274.2
$ws.Range("L137").Value = 16988.334
$ws.Range("M137").Value = -10174.2
$ws.Range("N137").Value = -27188.334

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1476.84
$ws.Range("I102").Value = 1409.2084
$ws.Range("J102").Value = 3100
$ws.Range("K102").Value = 1409.2084
$ws.Range("L102").Value = 3100
$ws.Range("M102").Value = 212.7916
$ws.Range("N102").Value = -6344

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3715
$ws.Range("I16").Value = 3426
$ws.Range("J16").Value = 5449
$ws.Range("K16").Value = 3426
$ws.Range("L16").Value = 5449
$ws.Range("M16").Value = -3256
$ws.Range("N16").Value = -5789
$ws.Range("H25").Value = 12000
$ws.Range("J25").Value = 12000
$ws.Range("L25").Value = 12000
$ws.Range("N25").Value = -12460
$ws.Range("H40").Value = 2484.913
$ws.Range("I40").Value = 2088.8333
$ws.Range("K40").Value = 2088.8333
$ws.Range("M40").Value = -1952.8333
$ws.Range("H46").Value = 1755.5714
$ws.Range("I46").Value = 1848.1666
$ws.Range("K46").Value = 1848.1666
$ws.Range("M46").Value = -1660.1666
$ws.Range("H93").Value = 1499.4286
$ws.Range("I93").Value = 1391
$ws.Range("J93").Value = 2150
$ws.Range("K93").Value = 1391
$ws.Range("L93").Value = 2150
$ws.Range("M93").Value = -143
$ws.Range("N93").Value = -4646
$ws.Range("H122").Value = 8240.526
$ws.Range("I122").Value = 8622
$ws.Range("K122").Value = 25866
$ws.Range("M122").Value = -23416
$ws.Range("H133").Value = 89993
$ws.Range("J133").Value = 89993
$ws.Range("L133").Value = 89993
$ws.Range("N133").Value = -95053
$ws.Range("H136").Value = 6893.75
$ws.Range("I136").Value = 5021.4287
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 15064.2861
$ws.Range("L136").Value = 60000
$ws.Range("M136").Value = -12514.2861
$ws.Range("N136").Value = -65100

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6125.3447
$ws.Range("I132").Value = 3806.111
$ws.Range("J132").Value = 9920.454
$ws.Range("K132").Value = 11418.333
$ws.Range("L132").Value = 29761.362
$ws.Range("M132").Value = -8888.332999999999
$ws.Range("N132").Value = -34821.362
$ws.Range("H135").Value = 75358
$ws.Range("J135").Value = 75358
$ws.Range("L135").Value = 75358
$ws.Range("N135").Value = -85498
$ws.Range("H136").Value = 4404.5356
$ws.Range("I136").Value = 4588.364
$ws.Range("K136").Value = 13765.092
$ws.Range("M136").Value = -11215.092
